$wb = $excel.ActiveWorkbook

# Add the new worksheet "Tabelle2" after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Tabelle2"

# Fill A1:E1 with "Test"
$newSheet.Range("A1:E1").Value = "Test"

# Make the new sheet the active / selected sheet
$newSheet.Select()
[void]$newSheet.Range("F1").Select()
